# Generate Report for Handback
#
# This script mutates the "localization-status.xlsx" handback report to:
#   1) Flip the Overview sheet's per-language status text from
#      "Ready for handoff" to "Handed back: in sync with en-US".
#   2) Record the handback timestamps for the zh-cn and de-de language
#      tables (Latest Handback DateTime column) - zh-cn finished syncing
#      at 2016-08-23 02:29:15, de-de at 2016-08-23 02:29:21.
#   3) Fill in the "Latest Target File" / "Latest Handback File" columns
#      (I/J) for both language tables, with I linked back to the source
#      markdown file on GitHub (mirroring the existing column A links).
#   4) Widen a handful of columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$zhHandbackDateTime = "2016-08-23 02:29:15"
$deHandbackDateTime = "2016-08-23 02:29:21"

$mdFile320 = "320f7643-be62-4404-86af-3eaf4f535905.md"
$mdFile8d7 = "8d755339-df3d-44fd-91d0-d5778c874b97.md"

$zhXlf320 = "320f7643-be62-4404-86af-3eaf4f535905.2342a7e23e195d5b1edbd169a3be67c8a0310178.zh-cn.xlf"
$zhXlf8d7 = "8d755339-df3d-44fd-91d0-d5778c874b97.5bfe23e8ba9b7dd23ba069607b0fdd075ce76249.zh-cn.xlf"
$deXlf320 = "320f7643-be62-4404-86af-3eaf4f535905.2342a7e23e195d5b1edbd169a3be67c8a0310178.de-de.xlf"
$deXlf8d7 = "8d755339-df3d-44fd-91d0-d5778c874b97.5bfe23e8ba9b7dd23ba069607b0fdd075ce76249.de-de.xlf"

$url320 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/215c487180b41a178aaf9853314317b96cb118f0/e2e/320f7643-be62-4404-86af-3eaf4f535905.md"
$url8d7 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/215c487180b41a178aaf9853314317b96cb118f0/e2e/8d755339-df3d-44fd-91d0-d5778c874b97.md"

# --- 1) Overview sheet: status column (E/F) for both tracked files -------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- 2) zh-cn table: status + target/handback file + handback datetime --
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Range("I2").Value = $mdFile320
$wsZhCn.Range("J2").Value = $zhXlf320
$wsZhCn.Range("K2").Value = $zhHandbackDateTime

$wsZhCn.Range("I3").Value = $mdFile8d7
$wsZhCn.Range("J3").Value = $zhXlf8d7
$wsZhCn.Range("K3").Value = $zhHandbackDateTime

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $url320, "", "", $mdFile320)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $url8d7, "", "", $mdFile8d7)

# --- 3) de-de table: status + target/handback file + handback datetime --
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Range("I2").Value = $mdFile320
$wsDeDe.Range("J2").Value = $deXlf320
$wsDeDe.Range("K2").Value = $deHandbackDateTime

$wsDeDe.Range("I3").Value = $mdFile8d7
$wsDeDe.Range("J3").Value = $deXlf8d7
$wsDeDe.Range("K3").Value = $deHandbackDateTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $url320, "", "", $mdFile320)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $url8d7, "", "", $mdFile8d7)

# --- 4) Column widths (widened to fit the newly populated content) ------
# Excel quantizes ColumnWidth to 1/6-character steps on save, so feed it
# already-quantized inputs (target - 5/6) to land as close as possible to
# the desired stored width.
$wide30 = 29.166666666666668   # -> stored width 30   (was ~17.22)
$wide40 = 39.166666666666664   # -> stored width 40   (was ~18.65 / 21.71)

$wsOverview.Columns.Item(5).ColumnWidth = $wide30   # column E
$wsOverview.Columns.Item(6).ColumnWidth = $wide30   # column F

$wsZhCn.Columns.Item(3).ColumnWidth = $wide30        # column C (Status)
$wsZhCn.Columns.Item(9).ColumnWidth = $wide40        # column I (Latest Target File)
$wsZhCn.Columns.Item(10).ColumnWidth = $wide40       # column J (Latest Handback File)

$wsDeDe.Columns.Item(3).ColumnWidth = $wide30        # column C (Status)
$wsDeDe.Columns.Item(9).ColumnWidth = $wide40        # column I (Latest Target File)
$wsDeDe.Columns.Item(10).ColumnWidth = $wide40       # column J (Latest Handback File)
